# The "Lãnh đạo hiện tại" row described both the Chairman and the Rector
# in a single cell (row 4: A4/B4). Split it into two rows: one for the
# Chairman ("Chủ tịch hội đồng quản trị") kept on row 4, and a new row 5
# for the Rector ("Hiệu trưởng hiện tại"). Every row below shifts down by
# one, including the "Hợp tác quốc tế" row whose column C carries a
# hyperlink (moves from C14 to C15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after row 4 so existing row 4 can be edited in
# place and the split-off Rector info lands on the freshly inserted row.
$ws.Rows(5).EntireRow.Insert()

# Row 4: now holds only the Chairman info.
$ws.Range("A4").Value = "Chủ tịch hội đồng quản trị"
$ws.Range("B4").Value = "Chủ tịch HĐQT: PGS. Lê Công Cơ. "

# Row 5 (newly inserted): holds the Rector info that used to be appended to B4.
$ws.Range("A5").Value = "Hiệu trưởng hiện tại"
$ws.Range("B5").Value = "Hiệu trưởng: con trai của thầy Lê Công Cơ (TS. Lê Nguyên Bảo)"

# The row insert shifted the "Hợp tác quốc tế" international-partners row
# (with its rich-text / hyperlinked C column) from row 14 down to row 15.
# The engine does not automatically re-point the worksheet hyperlink to
# follow the shifted cell, so re-create it at the new location pointing at
# the same external address, preserving the display text.
$displayText = "Carnegie Mellon ( CMU): https://duytan.edu.vn/tuyen-sinh/Page/EnrollArticleViewDetail.aspx?id=461`nPen State( PSU): https://duytan.edu.vn/tuyen-sinh/Page/EnrollArticleViewDetail.aspx?id=462 `nCalState Fullerton ( CSU): https://duytan.edu.vn/tuyen-sinh/Page/EnrollArticleViewDetail.aspx?id=463`nPurdue Northwest ( PNU): https://duytan.edu.vn/tuyen-sinh/Page/EnrollArticleViewDetail.aspx?id=819`nTroy: https://duytan.edu.vn/news/NewsDetail.aspx?id=4994&pid=2064&lang=vi-VN`nKeuka: https://duytan.edu.vn/news/NewsDetail.aspx?id=4994&pid=2064&lang=vi-VN`nDu học 3+1: https://duytan.edu.vn/tuyen-sinh/Page/EducationDetail.aspx?id=71`nDu học 2+2: https://duytan.edu.vn/tuyen-sinh/Page/EducationDetail.aspx?id=81`nDu học 1+1+2: https://duytan.edu.vn/tuyen-sinh/Page/EducationDetail.aspx?id=62"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C15"), "https://duytan.edu.vn/tuyen-sinh/Page/EnrollArticleViewDetail.aspx?id=461", "", "", $displayText)
